$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.992188944464516
$ws.Cells.Item(2, 3).Value = 0.04217500932493579
$ws.Cells.Item(2, 4).Value = 0.09010727540739083
$ws.Cells.Item(2, 6).Value = 1.82189692221371
$ws.Cells.Item(2, 7).Value = 0.002526713677461422
$ws.Cells.Item(2, 9).Value = 1.488173195535339
$ws.Cells.Item(2, 11).Value = 0.7524519804105125
$ws.Cells.Item(2, 12).Value = 0.305419044121507
$ws.Cells.Item(2, 14).Value = 2.539220195181542
$ws.Cells.Item(3, 2).Value = 0.9459752841745512
$ws.Cells.Item(3, 3).Value = 0.03700281875606493
$ws.Cells.Item(3, 4).Value = 0.09026088554658074
$ws.Cells.Item(3, 6).Value = 1.807764980417915
$ws.Cells.Item(3, 7).Value = 0.002530571122601493
$ws.Cells.Item(3, 9).Value = 1.485483617240639
$ws.Cells.Item(3, 11).Value = 0.7051070124853709
$ws.Cells.Item(3, 12).Value = 0.2946358510807983
$ws.Cells.Item(3, 14).Value = 2.553437943432222
$ws.Cells.Item(4, 2).Value = 0.9181663518681376
$ws.Cells.Item(4, 3).Value = 0.03381546066887609
$ws.Cells.Item(4, 4).Value = 0.09035624849220625
$ws.Cells.Item(4, 6).Value = 1.799999656348035
$ws.Cells.Item(4, 7).Value = 0.002533065751526336
$ws.Cells.Item(4, 9).Value = 1.484479200848874
$ws.Cells.Item(4, 11).Value = 0.6764588377691609
$ws.Cells.Item(4, 12).Value = 0.2881963055578751
$ws.Cells.Item(4, 14).Value = 2.562874846457007
$ws.Cells.Item(5, 2).Value = 0.9069765079668173
$ws.Cells.Item(5, 3).Value = 0.03251359684652755
$ws.Cells.Item(5, 4).Value = 0.09039537693603172
$ws.Cells.Item(5, 6).Value = 1.797064258851918
$ws.Cells.Item(5, 7).Value = 0.002534114153021405
$ws.Cells.Item(5, 9).Value = 1.484232477529922
$ws.Cells.Item(5, 11).Value = 0.6648904733922905
$ws.Cells.Item(5, 12).Value = 0.2856176885319144
$ws.Cells.Item(5, 14).Value = 2.566898127129662
$ws.Cells.Item(6, 2).Value = 0.9051270535697995
$ws.Cells.Item(6, 3).Value = 0.03229724076604157
$ws.Cells.Item(6, 4).Value = 0.09040189045634861
$ws.Cells.Item(6, 6).Value = 1.796590664856154
$ws.Cells.Item(6, 7).Value = 0.002534290164161753
$ws.Cells.Item(6, 9).Value = 1.484201325931821
$ws.Cells.Item(6, 11).Value = 0.6629759554576822
$ws.Cells.Item(6, 12).Value = 0.2851922605629795
$ws.Cells.Item(6, 14).Value = 2.567576916019789
$ws.Cells.Item(7, 2).Value = 0.9180148645377244
$ws.Cells.Item(7, 3).Value = 0.03379791547409638
$ws.Cells.Item(7, 4).Value = 0.09035677510484064
$ws.Cells.Item(7, 6).Value = 1.799959141508012
$ws.Cells.Item(7, 7).Value = 0.002533079761581106
$ws.Cells.Item(7, 9).Value = 1.48447521530192
$ws.Cells.Item(7, 11).Value = 0.6763023936661057
$ws.Cells.Item(7, 12).Value = 0.2881613450724387
$ws.Cells.Item(7, 14).Value = 2.562928386677775
$ws.Cells.Item(8, 2).Value = 0.9761369422066366
$ws.Cells.Item(8, 3).Value = 0.04039400480102984
$ws.Cells.Item(8, 4).Value = 0.09016002572307968
$ws.Cells.Item(8, 6).Value = 1.816834847327172
$ws.Cells.Item(8, 7).Value = 0.002528017604499625
$ws.Cells.Item(8, 9).Value = 1.487111475325761
$ws.Cells.Item(8, 11).Value = 0.7360398240411996
$ws.Cells.Item(8, 12).Value = 0.3016633141103426
$ws.Cells.Item(8, 14).Value = 2.543975614402854
$ws.Cells.Item(9, 2).Value = 1.094612284411511
$ws.Cells.Item(9, 3).Value = 0.0532404949301224
$ws.Cells.Item(9, 4).Value = 0.08978229980558972
$ws.Cells.Item(9, 6).Value = 1.85717680244089
$ws.Cells.Item(9, 7).Value = 0.002519087018186185
$ws.Cells.Item(9, 9).Value = 1.497421653740439
$ws.Cells.Item(9, 11).Value = 0.8565421567702174
$ws.Cells.Item(9, 12).Value = 0.3295842577677632
$ws.Cells.Item(9, 14).Value = 2.512427099460638
$ws.Cells.Item(10, 2).Value = 1.184414206860822
$ws.Cells.Item(10, 3).Value = 0.06263088563437691
$ws.Cells.Item(10, 4).Value = 0.08950944299350461
$ws.Cells.Item(10, 6).Value = 1.891260935272655
$ws.Cells.Item(10, 7).Value = 0.002513126646994476
$ws.Cells.Item(10, 9).Value = 1.50814270694643
$ws.Cells.Item(10, 11).Value = 0.9471466927828089
$ws.Cells.Item(10, 12).Value = 0.3509866276233282
$ws.Cells.Item(10, 14).Value = 2.492682825267721
$ws.Cells.Item(11, 2).Value = 1.225870479762136
$ws.Cells.Item(11, 3).Value = 0.06689374879958621
$ws.Cells.Item(11, 4).Value = 0.08938626800408578
$ws.Cells.Item(11, 6).Value = 1.907737862487352
$ws.Cells.Item(11, 7).Value = 0.002510544236078438
$ws.Cells.Item(11, 9).Value = 1.513706056894392
$ws.Cells.Item(11, 11).Value = 0.9888208439592745
$ws.Cells.Item(11, 12).Value = 0.3609181878510412
$ws.Cells.Item(11, 14).Value = 2.484448567300007
$ws.Cells.Item(12, 2).Value = 1.241655993069514
$ws.Cells.Item(12, 3).Value = 0.06850679765432233
$ws.Cells.Item(12, 4).Value = 0.08933975746369249
$ws.Cells.Item(12, 6).Value = 1.914117381144848
$ws.Cells.Item(12, 7).Value = 0.002509584788745479
$ws.Cells.Item(12, 9).Value = 1.515911629067332
$ws.Cells.Item(12, 11).Value = 1.004667845830795
$ws.Cells.Item(12, 12).Value = 0.3647072390475898
$ws.Cells.Item(12, 14).Value = 2.481438157658516
$ws.Cells.Item(13, 2).Value = 1.238252435428137
$ws.Cells.Item(13, 3).Value = 0.0681594511774648
$ws.Cells.Item(13, 4).Value = 0.08934976846916065
$ws.Cells.Item(13, 6).Value = 1.912737203321711
$ws.Cells.Item(13, 7).Value = 0.002509790603501121
$ws.Cells.Item(13, 9).Value = 1.515432220497061
$ws.Cells.Item(13, 11).Value = 1.001251977811222
$ws.Cells.Item(13, 12).Value = 0.3638899444213592
$ws.Cells.Item(13, 14).Value = 2.482081708360212
$ws.Cells.Item(14, 2).Value = 1.227167421744127
$ws.Cells.Item(14, 3).Value = 0.06702647903946968
$ws.Cells.Item(14, 4).Value = 0.08938243890354158
$ws.Cells.Item(14, 6).Value = 1.908259900457111
$ws.Cells.Item(14, 7).Value = 0.002510464932203066
$ws.Cells.Item(14, 9).Value = 1.513885528543156
$ws.Cells.Item(14, 11).Value = 0.9901232649726239
$ws.Cells.Item(14, 12).Value = 0.3612293501009702
$ws.Cells.Item(14, 14).Value = 2.484198738760554
$ws.Cells.Item(15, 2).Value = 1.220388848802543
$ws.Cells.Item(15, 3).Value = 0.06633234608997896
$ws.Cells.Item(15, 4).Value = 0.08940246773358673
$ws.Cells.Item(15, 6).Value = 1.905535672848174
$ws.Cells.Item(15, 7).Value = 0.002510880378863361
$ws.Cells.Item(15, 9).Value = 1.51295101348822
$ws.Cells.Item(15, 11).Value = 0.9833151925676589
$ws.Cells.Item(15, 12).Value = 0.3596033318581533
$ws.Cells.Item(15, 14).Value = 2.485509517014393
$ws.Cells.Item(16, 2).Value = 1.181717108924659
$ws.Cells.Item(16, 3).Value = 0.06235212501913168
$ws.Cells.Item(16, 4).Value = 0.08951751157646015
$ws.Cells.Item(16, 6).Value = 1.890203712366684
$ws.Cells.Item(16, 7).Value = 0.002513298001143433
$ws.Cells.Item(16, 9).Value = 1.507792952848305
$ws.Cells.Item(16, 11).Value = 0.9444324065082697
$ws.Cells.Item(16, 12).Value = 0.3503415170329731
$ws.Cells.Item(16, 14).Value = 2.493236015130577
$ws.Cells.Item(17, 2).Value = 1.158148137567821
$ws.Cells.Item(17, 3).Value = 0.0599081661921872
$ws.Cells.Item(17, 4).Value = 0.08958832788564308
$ws.Cells.Item(17, 6).Value = 1.881047191982987
$ws.Cells.Item(17, 7).Value = 0.00251481410585916
$ws.Cells.Item(17, 9).Value = 1.504804541090806
$ws.Cells.Item(17, 11).Value = 0.9206963652341358
$ws.Cells.Item(17, 12).Value = 0.3447098192182239
$ws.Cells.Item(17, 14).Value = 2.498167592194818
$ws.Cells.Item(18, 2).Value = 1.144648830909603
$ws.Cells.Item(18, 3).Value = 0.05850162631611511
$ws.Cells.Item(18, 4).Value = 0.08962914905414365
$ws.Cells.Item(18, 6).Value = 1.875872063265376
$ws.Cells.Item(18, 7).Value = 0.002515698276334815
$ws.Cells.Item(18, 9).Value = 1.503150272400781
$ws.Cells.Item(18, 11).Value = 0.9070870967299243
$ws.Cells.Item(18, 12).Value = 0.3414890141261964
$ws.Cells.Item(18, 14).Value = 2.501074449901466
$ws.Cells.Item(19, 2).Value = 1.140087979593147
$ws.Cells.Item(19, 3).Value = 0.0580252497569802
$ws.Cells.Item(19, 4).Value = 0.0896429858771004
$ws.Cells.Item(19, 6).Value = 1.874135552573179
$ws.Cells.Item(19, 7).Value = 0.002515999729939766
$ws.Cells.Item(19, 9).Value = 1.502601253315021
$ws.Cells.Item(19, 11).Value = 0.9024866277077308
$ws.Cells.Item(19, 12).Value = 0.3404016627082171
$ws.Cells.Item(19, 14).Value = 2.502070736755925
$ws.Cells.Item(20, 2).Value = 1.160651201980329
$ws.Cells.Item(20, 3).Value = 0.06016841614332691
$ws.Cells.Item(20, 4).Value = 0.08958078012872406
$ws.Cells.Item(20, 6).Value = 1.882012451255008
$ws.Cells.Item(20, 7).Value = 0.002514651457786446
$ws.Cells.Item(20, 9).Value = 1.505115976950748
$ws.Cells.Item(20, 11).Value = 0.9232186486515275
$ws.Cells.Item(20, 12).Value = 0.3453074183825606
$ws.Cells.Item(20, 14).Value = 2.497635335257542
$ws.Cells.Item(21, 2).Value = 1.230421000413628
$ws.Cells.Item(21, 3).Value = 0.06735929255357576
$ws.Cells.Item(21, 4).Value = 0.08937283921665795
$ws.Cells.Item(21, 6).Value = 1.909571189067492
$ws.Cells.Item(21, 7).Value = 0.002510266365339251
$ws.Cells.Item(21, 9).Value = 1.514337145311089
$ws.Cells.Item(21, 11).Value = 0.9933902483120107
$ws.Cells.Item(21, 12).Value = 0.3620100656536493
$ws.Cells.Item(21, 14).Value = 2.483573990098833
$ws.Cells.Item(22, 2).Value = 1.276526220799326
$ws.Cells.Item(22, 3).Value = 0.0720519560378392
$ws.Cells.Item(22, 4).Value = 0.08923771254456447
$ws.Cells.Item(22, 6).Value = 1.928398963173464
$ws.Cells.Item(22, 7).Value = 0.002507507986946044
$ws.Cells.Item(22, 9).Value = 1.520939977322428
$ws.Cells.Item(22, 11).Value = 1.039635636203798
$ws.Cells.Item(22, 12).Value = 0.3730905274986753
$ws.Cells.Item(22, 14).Value = 2.475012104309059
$ws.Cells.Item(23, 2).Value = 1.251872661207699
$ws.Cells.Item(23, 3).Value = 0.06954800925552718
$ws.Cells.Item(23, 4).Value = 0.08930976234547217
$ws.Cells.Item(23, 6).Value = 1.91827540885663
$ws.Cells.Item(23, 7).Value = 0.002508970376780218
$ws.Cells.Item(23, 9).Value = 1.517363138066443
$ws.Cells.Item(23, 11).Value = 1.014918418113496
$ws.Cells.Item(23, 12).Value = 0.3671616181005106
$ws.Cells.Item(23, 14).Value = 2.479524202014588
$ws.Cells.Item(24, 2).Value = 1.159519408289952
$ws.Cells.Item(24, 3).Value = 0.06005076174361079
$ws.Cells.Item(24, 4).Value = 0.08958419213222513
$ws.Cells.Item(24, 6).Value = 1.881575780079899
$ws.Cells.Item(24, 7).Value = 0.002514724951743901
$ws.Cells.Item(24, 9).Value = 1.504974978037637
$ws.Cells.Item(24, 11).Value = 0.9220782093962328
$ws.Cells.Item(24, 12).Value = 0.3450371910651313
$ws.Cells.Item(24, 14).Value = 2.497875745493502
$ws.Cells.Item(25, 2).Value = 1.062078383589267
$ws.Cells.Item(25, 3).Value = 0.04977402834731492
$ws.Cells.Item(25, 4).Value = 0.08988364718895347
$ws.Cells.Item(25, 6).Value = 1.845484548650177
$ws.Cells.Item(25, 7).Value = 0.002521396986247832
$ws.Cells.Item(25, 9).Value = 1.494080884048202
$ws.Cells.Item(25, 11).Value = 0.8235811115017952
$ws.Cells.Item(25, 12).Value = 0.3218754488928965
$ws.Cells.Item(25, 14).Value = 2.520359401337259
